$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-08 -> 2023-09-09, serial 45177 -> 45178) for every data row
# (rows 2 through 395).
$ws.Range("C2:C395").Value = 45178
